$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999758567748,
    0.99756036203660892,
    0.99316843940224353,
    0.99580294866853747,
    0.98527854394369174,
    0.96022807009011935,
    0.9534656557579293,
    0.94433386252873452,
    0.93258589298702432,
    0.92175550349575086,
    0.92015911593213695,
    0.91737892860298698,
    0.90609059921803725,
    0.90192288999694736,
    0.89933120233880359,
    0.89682457632695911,
    0.89311648354843587,
    0.89200754169284968,
    0.99491018798415332,
    0.9736359606153846,
    0.96791585890027398,
    0.95870708233443258,
    0.98613416662613818,
    0.97311383839742271,
    0.96665692335182696,
    0.95870259354304954,
    0.9568270272396584,
    0.94917293658177693,
    0.94419134295953433,
    0.94275457799895168,
    0.93510179695473372,
    0.93342246139257157,
    0.93290243487163993
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
